# Insert a new row at position 36, shifting existing rows 36..113 down to 37..114
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(36).Insert()

# Fill the new row 36 with data (mostly a duplicate of the original row 36,
# now shifted to row 37, but with Fecha/Volumen/Precio/Origen/Precio-Kg updated)
$ws.Cells.Item(36, 1).Value = 8
$ws.Cells.Item(36, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(36, 3).Value = "Coquimbo"
$ws.Cells.Item(36, 4).Value = Get-Date -Year 2022 -Month 4 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(36, 5).Value = 4
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100109
$ws.Cells.Item(36, 8).Value = "Uva"
$ws.Cells.Item(36, 9).Value = 100109001
$ws.Cells.Item(36, 10).Value = "Uva"
$ws.Cells.Item(36, 11).Value = "Red Globe"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 400
$ws.Cells.Item(36, 14).Value = 9500
$ws.Cells.Item(36, 15).Value = 10000
$ws.Cells.Item(36, 16).Value = 9750
$ws.Cells.Item(36, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(36, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(36, 19).Value = 542
$ws.Cells.Item(36, 20).Value = 18
